$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra trial rows (rows 6-10), reducing the table from 9 trials to 4
$ws.Range("A6:C10").ClearContents()

# Keep 2 trials for letters (Phonemic) and 2 trials for categories (Semantic)
$ws.Cells.Item(2,1).Value() = "Words that start with A"
$ws.Cells.Item(2,2).Value() = "Phonemic"
$ws.Cells.Item(2,3).Value() = "j"

$ws.Cells.Item(3,1).Value() = "Words that start with S"
$ws.Cells.Item(3,2).Value() = "Phonemic"
$ws.Cells.Item(3,3).Value() = "p"

$ws.Cells.Item(4,1).Value() = "Animals"
$ws.Cells.Item(4,2).Value() = "Semantic"
$ws.Cells.Item(4,3).Value() = "j"

$ws.Cells.Item(5,1).Value() = "Occupations"
$ws.Cells.Item(5,2).Value() = "Semantic"
$ws.Cells.Item(5,3).Value() = "p"

# Update selection to match saved view state
$ws.Range("A9").Select()
